# Auto-applies the row-level numeric updates from the Kujata_Profits.xlsx diff
# across the ALC/ARM/BSM/CUL/GSM/LTW/WVR sheets. Values are written with
# Range.Value; cells the diff deletes are cleared with ClearContents().
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 3875.25
$ws.Range("I18").Value = 5500.5
$ws.Range("J18").Value = 2250
$ws.Range("K18").Value = 5500.5
$ws.Range("L18").Value = 2250
$ws.Range("M18").Value = -5216.5
$ws.Range("N18").Value = -2818

$ws.Range("H43").Value = 27778028
$ws.Range("J43").Value = 55555556
$ws.Range("L43").Value = 55555556
$ws.Range("N43").Value = -55555694

$ws.Range("H74").Value = 3575
$ws.Range("I74").Value = 3600
$ws.Range("J74").Value = 3500
$ws.Range("K74").Value = 3600
$ws.Range("L74").Value = 3500
$ws.Range("M74").Value = -2664
$ws.Range("N74").Value = -5372

$ws.Range("H77").Value = 3575
$ws.Range("I77").Value = 3600
$ws.Range("J77").Value = 3500
$ws.Range("K77").Value = 18000
$ws.Range("L77").Value = 17500
$ws.Range("M77").Value = -13320
$ws.Range("N77").Value = -26860

$ws.Range("H112").Value = 5522.25
$ws.Range("J112").Value = 6996.5
$ws.Range("L112").Value = 20989.5
$ws.Range("N112").Value = -23205.5

$ws.Range("H116").Value = 3419.9412
$ws.Range("I116").Value = 2753.7
$ws.Range("K116").Value = 2753.7
$ws.Range("M116").Value = 688.3000000000002

$ws.Range("H129").Value = 910.5714
$ws.Range("J129").Value = 913.9167
$ws.Range("L129").Value = 2741.7501
$ws.Range("N129").Value = -12741.7501

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4438.054
$ws.Range("I32").Value = 5070.567
$ws.Range("J32").Value = 1727.2858
$ws.Range("K32").Value = 5070.567
$ws.Range("L32").Value = 1727.2858
$ws.Range("M32").Value = -4783.567
$ws.Range("N32").Value = -2301.2858

$ws.Range("H74").Value = 1115.56
$ws.Range("I74").Value = 932.3684
$ws.Range("J74").Value = 1695.6666
$ws.Range("K74").Value = 932.3684
$ws.Range("L74").Value = 1695.6666
$ws.Range("M74").Value = -58.36839999999995
$ws.Range("N74").Value = -3443.6666

$ws.Range("H77").Value = 1115.56
$ws.Range("I77").Value = 932.3684
$ws.Range("J77").Value = 1695.6666
$ws.Range("K77").Value = 4661.842
$ws.Range("L77").Value = 8478.333000000001
$ws.Range("M77").Value = -293.8419999999996
$ws.Range("N77").Value = -17214.333

$ws.Range("H132").Value = 2824.325
$ws.Range("I132").Value = 2514.516
$ws.Range("J132").Value = 3891.4443
$ws.Range("K132").Value = 7543.548000000001
$ws.Range("L132").Value = 11674.3329
$ws.Range("M132").Value = -5013.548000000001
$ws.Range("N132").Value = -16734.3329

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws.Range("H94").Value = 12500789
$ws.Range("J94").Value = 1378
$ws.Range("L94").Value = 1378
$ws.Range("N94").Value = -2280

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 7115.3076
$ws.Range("I56").Value = 7115.3076
$ws.Range("K56").Value = 7115.3076
$ws.Range("M56").Value = -6585.3076

$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()

$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()

$ws.Range("H107").Value = 6339.1113
$ws.Range("I107").Value = 621.7143
$ws.Range("J107").Value = 9977.454
$ws.Range("K107").Value = 1865.1429
$ws.Range("L107").Value = 29932.362
$ws.Range("M107").Value = 54.85710000000017
$ws.Range("N107").Value = -33772.362

$ws.Range("H112").Value = 4158.6665
$ws.Range("J112").Value = 6166.6665
$ws.Range("L112").Value = 18499.9995
$ws.Range("N112").Value = -20715.9995

$ws.Range("H131").Value = 22223740
$ws.Range("J131").Value = 1590.8334
$ws.Range("L131").Value = 4772.5002
$ws.Range("N131").Value = -14852.5002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 18004398
$ws.Range("I70").Value = 20836922
$ws.Range("J70").Value = 15389760
$ws.Range("K70").Value = 20836922
$ws.Range("L70").Value = 15389760
$ws.Range("M70").Value = -20836652
$ws.Range("N70").Value = -15390300

$ws.Range("H73").Value = 18004398
$ws.Range("I73").Value = 20836922
$ws.Range("J73").Value = 15389760
$ws.Range("K73").Value = 20836922
$ws.Range("L73").Value = 15389760
$ws.Range("M73").Value = -20835986
$ws.Range("N73").Value = -15391632

$ws.Range("H80").Value = 5046
$ws.Range("I80").Value = 3899.875
$ws.Range("J80").Value = 6879.8
$ws.Range("K80").Value = 3899.875
$ws.Range("L80").Value = 6879.8
$ws.Range("M80").Value = -2901.875
$ws.Range("N80").Value = -8875.799999999999

$ws.Range("H83").Value = 5046
$ws.Range("I83").Value = 3899.875
$ws.Range("J83").Value = 6879.8
$ws.Range("K83").Value = 19499.375
$ws.Range("L83").Value = 34399
$ws.Range("M83").Value = -14507.375
$ws.Range("N83").Value = -44383

$ws.Range("H132").Value = 3414.353
$ws.Range("I132").Value = 2928.8333
$ws.Range("J132").Value = 4579.6
$ws.Range("K132").Value = 8786.499899999999
$ws.Range("L132").Value = 13738.8
$ws.Range("M132").Value = -6256.499899999999
$ws.Range("N132").Value = -18798.8

$ws.Range("H138").Value = 40109.668
$ws.Range("J138").Value = 40109.668
$ws.Range("L138").Value = 40109.668
$ws.Range("N138").Value = -50389.668

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1600
$ws.Range("J46").Value = 2100
$ws.Range("L46").Value = 2100
$ws.Range("N46").Value = -2476

$ws.Range("H100").Value = 2135
$ws.Range("I100").Value = 2001.5
$ws.Range("J100").Value = 2402
$ws.Range("K100").Value = 2001.5
$ws.Range("L100").Value = 2402
$ws.Range("M100").Value = -1460.5
$ws.Range("N100").Value = -3484

$ws.Range("H122").Value = 56673668
$ws.Range("J122").Value = 10000
$ws.Range("L122").Value = 30000
$ws.Range("N122").Value = -34900

$ws.Range("H132").Value = 52585.4
$ws.Range("I132").Value = 2189.7
$ws.Range("K132").Value = 6569.099999999999
$ws.Range("M132").Value = -4039.099999999999

$ws.Range("H136").Value = 5552.1665
$ws.Range("I136").Value = 7752.933
$ws.Range("K136").Value = 23258.799
$ws.Range("M136").Value = -20708.799

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 11767.111
$ws.Range("I132").Value = 16581.4
$ws.Range("J132").Value = 5749.25
$ws.Range("K132").Value = 49744.2
$ws.Range("L132").Value = 17247.75
$ws.Range("M132").Value = -22307.75
$ws.Range("N132").Value = -22307.75
